$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New feature/target pairs to append after the existing data (rows 2-220 already present, last row is 220)
$newRows = @(
    @("Quero agendar uma consulta", 0),
    @("Eu gostaria de agendar uma consulta para meu filho", 0),
    @("Limpeza , queria para o dia 10 de março", 0),
    @("Quero agendar", 0),
    @("Quero remarcar a minha consulta", 1),
    @("Eu queria remarcar a consulta do meu filho", 1),
    @("Quero remarcar", 1)
)

$startRow = 221
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}

# Add 15 trailing empty rows (228-242) with a distinct style (new font, no border)
$emptyStart = 228
$emptyEnd = 242
for ($r = $emptyStart; $r -le $emptyEnd; $r++) {
    $ws.Cells.Item($r, 1).Value = $null
    $ws.Cells.Item($r, 2).Value = $null
}
